$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2096.6
$ws.Range("I40").Value = 1782.8334
$ws.Range("J40").Value = 2567.25
$ws.Range("K40").Value = 1782.8334
$ws.Range("L40").Value = 2567.25
$ws.Range("M40").Value = -1607.8334
$ws.Range("N40").Value = -2917.25

$ws.Range("H58").Value = 3598.2
$ws.Range("J58").Value = 4845.2856
$ws.Range("L58").Value = 14535.8568
$ws.Range("N58").Value = -14835.8568

$ws.Range("H62").Value = 2747.8
$ws.Range("J62").Value = 1616.6666
$ws.Range("L62").Value = 1616.6666
$ws.Range("N62").Value = -2864.6666

$ws.Range("H65").Value = 2747.8
$ws.Range("J65").Value = 1616.6666
$ws.Range("L65").Value = 8083.333000000001
$ws.Range("N65").Value = -14323.333

$ws.Range("H86").Value = 4789.4707
$ws.Range("I86").Value = 2274.1
$ws.Range("J86").Value = 8382.857
$ws.Range("K86").Value = 2274.1
$ws.Range("L86").Value = 8382.857
$ws.Range("M86").Value = -1151.1
$ws.Range("N86").Value = -10628.857

$ws.Range("H89").Value = 4789.4707
$ws.Range("I89").Value = 2274.1
$ws.Range("J89").Value = 8382.857
$ws.Range("K89").Value = 11370.5
$ws.Range("L89").Value = 41914.285
$ws.Range("M89").Value = -5754.5
$ws.Range("N89").Value = -53146.285

$ws.Range("H97").Value = 2372.389
$ws.Range("J97").Value = 2372.389
$ws.Range("L97").Value = 7117.167
$ws.Range("N97").Value = -8109.167

$ws.Range("H111").Value = 2847.0454
$ws.Range("I111").Value = 1330.5454
$ws.Range("J111").Value = 4363.5454
$ws.Range("K111").Value = 3991.6362
$ws.Range("L111").Value = 13090.6362
$ws.Range("M111").Value = -924.6361999999999
$ws.Range("N111").Value = -19224.6362

$ws.Range("H112").Value = 2141.7
$ws.Range("I112").Value = 685.5714
$ws.Range("K112").Value = 2056.7142
$ws.Range("M112").Value = -948.7142000000003

$ws.Range("H138").Value = 2116.41
$ws.Range("I138").Value = 1589.9
$ws.Range("J138").Value = 2174.9111
$ws.Range("K138").Value = 4769.700000000001
$ws.Range("L138").Value = 6524.7333
$ws.Range("M138").Value = 370.2999999999993
$ws.Range("N138").Value = -16804.7333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 968.94446
$ws.Range("I74").Value = 1032.8462
$ws.Range("J74").Value = 802.8
$ws.Range("K74").Value = 1032.8462
$ws.Range("L74").Value = 802.8
$ws.Range("M74").Value = -158.8462
$ws.Range("N74").Value = -2550.8

$ws.Range("H77").Value = 968.94446
$ws.Range("I77").Value = 1032.8462
$ws.Range("J77").Value = 802.8
$ws.Range("K77").Value = 5164.231
$ws.Range("L77").Value = 4014
$ws.Range("M77").Value = -796.2309999999998
$ws.Range("N77").Value = -12750

$ws.Range("H108").Value = 32528
$ws.Range("J108").Value = 32528
$ws.Range("L108").Value = 32528
$ws.Range("N108").Value = -40208

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H139").Value = 50714
$ws.Range("J139").Value = 50714
$ws.Range("L139").Value = 50714
$ws.Range("N139").Value = -60994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 35715256
$ws.Range("I94").Value = 50000560
$ws.Range("K94").Value = 50000560
$ws.Range("M94").Value = -50000109

$ws.Range("H107").Value = 1642.6666
$ws.Range("I107").Value = 1382.7
$ws.Range("J107").Value = 2162.6
$ws.Range("K107").Value = 1382.7
$ws.Range("L107").Value = 2162.6
$ws.Range("M107").Value = 537.3
$ws.Range("N107").Value = -6002.6

$ws.Range("H134").Value = 5530.875
$ws.Range("I134").Value = 1216.45
$ws.Range("K134").Value = 3649.35
$ws.Range("M134").Value = -1114.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1304.8948
$ws.Range("I31").Value = 954.58826
$ws.Range("J31").Value = 1588.4762
$ws.Range("K31").Value = 954.58826
$ws.Range("L31").Value = 1588.4762
$ws.Range("M31").Value = -659.58826
$ws.Range("N31").Value = -2178.4762

$ws.Range("H34").Value = 1304.8948
$ws.Range("I34").Value = 954.58826
$ws.Range("J34").Value = 1588.4762
$ws.Range("K34").Value = 954.58826
$ws.Range("L34").Value = 1588.4762
$ws.Range("M34").Value = -752.58826
$ws.Range("N34").Value = -1992.4762

$ws.Range("H58").Value = 936.1875
$ws.Range("I58").Value = 1003.2381
$ws.Range("J58").Value = 808.1818
$ws.Range("K58").Value = 1003.2381
$ws.Range("L58").Value = 808.1818
$ws.Range("M58").Value = -800.2381
$ws.Range("N58").Value = -1214.1818

$ws.Range("H105").Value = 594.5
$ws.Range("I105").Value = 391.66666
$ws.Range("J105").Value = 746.625
$ws.Range("K105").Value = 391.66666
$ws.Range("L105").Value = 746.625
$ws.Range("M105").Value = 1355.33334
$ws.Range("N105").Value = -4240.625

$ws.Range("H133").Value = 61522.31
$ws.Range("J133").Value = 61522.31
$ws.Range("L133").Value = 61522.31
$ws.Range("N133").Value = -66582.31

$ws.Range("H136").Value = 936.1875
$ws.Range("I136").Value = 1003.2381
$ws.Range("J136").Value = 808.1818
$ws.Range("K136").Value = 3009.7143
$ws.Range("L136").Value = 2424.5454
$ws.Range("M136").Value = -459.7143000000001
$ws.Range("N136").Value = -7524.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1716.7273
$ws.Range("I99").Value = 380
$ws.Range("K99").Value = 1140
$ws.Range("M99").Value = 1106

$ws.Range("H113").Value = 556.8333
$ws.Range("J113").Value = 608.58826
$ws.Range("L113").Value = 1825.76478
$ws.Range("N113").Value = -6165.76478

$ws.Range("H136").Value = 1421.5714
$ws.Range("J136").Value = 4533
$ws.Range("L136").Value = 13599
$ws.Range("N136").Value = -23799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2333.3333
$ws.Range("I9").Value = 2333.3333
$ws.Range("K9").Value = 2333.3333
$ws.Range("M9").Value = -2163.3333

$ws.Range("H80").Value = 5550
$ws.Range("J80").Value = 6800
$ws.Range("L80").Value = 6800
$ws.Range("N80").Value = -8796

$ws.Range("H83").Value = 5550
$ws.Range("J83").Value = 6800
$ws.Range("L83").Value = 34000
$ws.Range("N83").Value = -43984

$ws.Range("H104").Value = 57333.332
$ws.Range("J104").Value = 57333.332
$ws.Range("L104").Value = 57333.332
$ws.Range("N104").Value = -64321.332

$ws.Range("H113").Value = 1807
$ws.Range("I113").Value = 1056.1111
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 1056.1111
$ws.Range("L113").Value = 2933.3333
$ws.Range("M113").Value = 1113.8889
$ws.Range("N113").Value = -7273.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 688.7222
$ws.Range("I22").Value = 319.8
$ws.Range("J22").Value = 1149.875
$ws.Range("K22").Value = 319.8
$ws.Range("L22").Value = 1149.875
$ws.Range("M22").Value = -24.80000000000001
$ws.Range("N22").Value = -1739.875

$ws.Range("H27").Value = 688.7222
$ws.Range("I27").Value = 319.8
$ws.Range("J27").Value = 1149.875
$ws.Range("K27").Value = 319.8
$ws.Range("L27").Value = 1149.875
$ws.Range("M27").Value = -212.8
$ws.Range("N27").Value = -1363.875

$ws.Range("H61").Value = 2575.3
$ws.Range("I61").Value = 2109.7
$ws.Range("J61").Value = 3040.9
$ws.Range("K61").Value = 2109.7
$ws.Range("L61").Value = 3040.9
$ws.Range("M61").Value = -1907.7
$ws.Range("N61").Value = -3444.9

$ws.Range("H100").Value = 2450.4
$ws.Range("J100").Value = 2643.4285
$ws.Range("L100").Value = 2643.4285
$ws.Range("N100").Value = -3725.4285

$ws.Range("H113").Value = 2575.3
$ws.Range("I113").Value = 2109.7
$ws.Range("J113").Value = 3040.9
$ws.Range("K113").Value = 2109.7
$ws.Range("L113").Value = 3040.9
$ws.Range("M113").Value = 60.30000000000018
$ws.Range("N113").Value = -7380.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 41464.5
$ws.Range("J139").Value = 48215
$ws.Range("L139").Value = 48215
$ws.Range("N139").Value = -58495
